$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.046.91'
$ws.Range("E2").Value = '  -1.73%  '
$ws.Range("D3").Value = '2.420.03'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.99'
$ws.Range("E5").Value = '  -2.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.02'
$ws.Range("E6").Value = '  -1.95%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.527'
$ws.Range("E8").Value = '  -0.58%  '
$ws.Range("D9").Value = '2.404.94'
$ws.Range("E9").Value = '  -1.41%  '
$ws.Range("E10").Value = '  -1.10%  '
$ws.Range("E11").Value = '  -0.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.08'
$ws.Range("E12").Value = '  -2.34%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.339'
$ws.Range("E13").Value = '  -1.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.19'
$ws.Range("E14").Value = '  -1.20%  '
$ws.Range("E15").Value = '  -2.81%  '
$ws.Range("D16").Value = '2.829.24'
$ws.Range("E16").Value = '  -2.07%  '
$ws.Range("D17").Value = '60.889.05'
$ws.Range("E17").Value = '  -1.81%  '
$ws.Range("D18").Value = '2.408.19'
$ws.Range("E18").Value = '  -1.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.73'
$ws.Range("E19").Value = '  +8.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.66'
$ws.Range("E20").Value = '  -1.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.55'
$ws.Range("E21").Value = '  -0.91%  '
$ws.Range("E22").Value = '  -1.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.09'
$ws.Range("E23").Value = '  +1.80%  '
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.85'
$ws.Range("E25").Value = '  -3.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '64.92'
$ws.Range("E26").Value = '  -1.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '584.92'
$ws.Range("E27").Value = '  -1.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.35'
$ws.Range("E28").Value = '  -9.23%  '
$ws.Range("D29").Value = '2.536.80'
$ws.Range("E29").Value = '  -1.36%  '
$ws.Range("D30").Value = '0.0₃0936'
$ws.Range("E30").Value = '  -3.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.90'
$ws.Range("E31").Value = '  -0.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.35'
$ws.Range("E32").Value = '  -4.47%  '
$ws.Range("E33").Value = '  -2.65%  '
$ws.Range("E34").Value = '  -2.50%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.41'
$ws.Range("E36").Value = '  -1.02%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.62'
$ws.Range("E37").Value = '  -5.23%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '151.37'
$ws.Range("E38").Value = '  -2.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.368'
$ws.Range("E39").Value = '  -1.87%  '
$ws.Range("E40").Value = '  -0.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.15'
$ws.Range("E41").Value = '  -2.43%  '
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("E43").Value = '  -1.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.14'
$ws.Range("E44").Value = '  -5.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.37'
$ws.Range("E45").Value = '  -5.42%  '
$ws.Range("D46").Value = '0.0₆0273'
$ws.Range("E46").Value = '  +1.81%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '142.58'
$ws.Range("E47").Value = '  +0.58%  '
$ws.Range("E48").Value = '  -3.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.586'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.49'
$ws.Range("E50").Value = '  -1.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0503'
$ws.Range("E51").Value = '  -3.26%  '
